$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.378347992897034
$ws.Range("B1").Value = 4.233298301696777
$ws.Range("C1").Value = 2.147314548492432
$ws.Range("D1").Value = 1.660015344619751
$ws.Range("E1").Value = 1.501710891723633
